# Update "Prix Spot" sheet: add a new day column BC (07-aug) with hourly prices
$wb = $excel.ActiveWorkbook

$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell BC1 - replicate the style (bold / bordered / centered) of the
# preceding header cell BB1 before writing the new header text.
$wsPrix.Range("BB1").Copy()
$wsPrix.Range("BC1").PasteSpecial(-4122)
$wsPrix.Range("BC1").Value = "07-aug"

# Hourly values for 07-aug in column BC (column 55), rows 2-25
$prixValues = @{
    2  = 88.42
    3  = 61.2
    4  = 44.63
    5  = 52.66
    6  = 46.6
    7  = 41.73
    8  = 42.97
    9  = 71.97
    10 = 77.02
    11 = 71.65000000000001
    12 = 22.23
    13 = 0
    14 = -0.01
    15 = -0.02
    16 = -0.07000000000000001
    17 = -0.01
    18 = 3.81
    19 = 48.97
    20 = 70
    21 = 101.53
    22 = 115.42
    23 = 110.74
    24 = 110
    25 = 101.28
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Cells.Item($row, 55).Value = $prixValues[$row]
}

# Update "Gaz" sheet: add new row 52 (2025-08-05)
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A52").NumberFormat = "@"
$wsGaz.Range("A52").Value = "2025-08-05"
$wsGaz.Range("A52").Style = "Normal"
$wsGaz.Range("B52").Value = 33.775

# Update "CO2" sheet: add new row 52 (2025-08-05) with an empty price value
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A52").NumberFormat = "@"
$wsCo2.Range("A52").Value = "2025-08-05"
$wsCo2.Range("A52").Style = "Normal"
$wsCo2.Range("B52").Value = "'"
$wsCo2.Range("B52").Style = "Normal"
